# M10 Froze Encoder 1234
# Update column C values on the active sheet to reflect new ASR results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 8
$ws.Range("C3").Value = 9
$ws.Range("C4").Value = 5
$ws.Range("C5").Value = 8
$ws.Range("C6").Value = 8
$ws.Range("C8").Value = 5
$ws.Range("C10").Value = 6
$ws.Range("C11").Value = 5
$ws.Range("C12").Value = 5
$ws.Range("C13").Value = 11
$ws.Range("C14").Value = 5
$ws.Range("C15").Value = 4
$ws.Range("C16").Value = 6
$ws.Range("C18").Value = 9
